# Apply cell-value updates from the "Updated cryptos list" refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds numeric-looking values that must stay plain text
# (inline/shared strings in the source workbook), so force text format before
# writing them and then drop back to the Normal style to avoid leaving a custom
# number-format behind on the cell.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "37.402.33"
$ws.Range("E2").Value = "  +0.15%  "
Set-TextValue "D3" "2.067.83"
$ws.Range("E3").Value = "  +0.48%  "
$ws.Range("E4").Value = "  +0.07%  "
Set-TextValue "D5" "235.27"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("E7").Value = "  +0.05%  "
Set-TextValue "D8" "57.60"
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("E9").Value = "  +3.62%  "
Set-TextValue "D10" "0.0773"
$ws.Range("E10").Value = "  +1.59%  "
Set-TextValue "D11" "0.103"
$ws.Range("E11").Value = "  +0.85%  "
Set-TextValue "D12" "2.373.72"
$ws.Range("E12").Value = "  +0.54%  "
Set-TextValue "D13" "14.47"
$ws.Range("E13").Value = "  -0.03%  "
Set-TextValue "D14" "20.82"
$ws.Range("E14").Value = "  -0.78%  "
Set-TextValue "D15" "0.785"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("E16").Value = "  +0.51%  "
Set-TextValue "D17" "2.068.58"
$ws.Range("E17").Value = "  +0.37%  "
Set-TextValue "D18" "37.358.31"
$ws.Range("E18").Value = "  -0.45%  "
Set-TextValue "D19" "6.34"
$ws.Range("E19").Value = "  +3.20%  "
Set-TextValue "D20" "69.73"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("E21").Value = "  +0.25%  "
Set-TextValue "D22" "226.90"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("E25").Value = "  -1.38%  "
Set-TextValue "D26" "167.24"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  -4.56%  "
Set-TextValue "D29" "19.12"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E31").Value = "  -1.09%  "
Set-TextValue "D32" "4.56"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("E33").Value = "  -0.57%  "
Set-TextValue "D34" "4.56"
$ws.Range("E34").Value = "  +1.70%  "
Set-TextValue "D35" "2.51"
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("E36").Value = "  +0.12%  "
Set-TextValue "D37" "3.32"
$ws.Range("E37").Value = "  -3.04%  "
$ws.Range("E38").Value = "  +0.15%  "
Set-TextValue "D39" "5.65"
$ws.Range("E39").Value = "  -3.76%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D40" "2.95"
$ws.Range("E40").Value = "  -0.66%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D41" "0.0965"
$ws.Range("E41").Value = "  -2.50%  "
Set-TextValue "D42" "98.32"
$ws.Range("E42").Value = "  +1.37%  "
Set-TextValue "D43" "1.482.57"
$ws.Range("E43").Value = "  +0.64%  "
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  -0.18%  "
Set-TextValue "D46" "4.05"
$ws.Range("E46").Value = "  -12.31%  "
Set-TextValue "D48" "15.35"
$ws.Range("E48").Value = "  -3.33%  "
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("E50").Value = "  +0.97%  "
Set-TextValue "D51" "2.260.08"
$ws.Range("E51").Value = "  +0.46%  "
